# 05.jedis操作redis集群.docx -- from wangyingchen at 20191225
#
# The original runs get re-split the way Word's background spell-checker
# splits them once it has re-scanned the paragraph (mixed CJK/Latin text
# triggers "spellStart"/"spellEnd" proofErr bracketing around the Latin
# words), and a new citation of the two source files backing the "【参考】"
# note is appended. A trailing empty paragraph is also added at the end
# of the document.
#
# Because the COM surface doesn't expose a way to poke a lone <w:proofErr/>
# in between two existing runs, each touched paragraph is rewritten in one
# shot with Range.InsertXML (it replaces the exact range's contents, so
# targeting Paragraphs(n).Range keeps the edit scoped to that paragraph).

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- Paragraph 1: "Jedis( JedisCluster )操作Redis集群redis-cluster" ---------
$p1 = $d.Paragraphs(1).Range
$xml1 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00055301" w:rsidRDefault="00B72A5A">' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Jedis</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">( </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>JedisCluster</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> )</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>操作</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Redis</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>集群</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>redis</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>-cluster</w:t></w:r>' + `
  '</w:p>'
$p1.InsertXML($xml1)

# --- Paragraph 2: "Java( JedisCluster )操作redis集群" -----------------------
$p2 = $d.Paragraphs(2).Range
$xml2 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00B72A5A" w:rsidRDefault="00B72A5A">' + `
  '<w:r><w:t xml:space="preserve">Java( </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>JedisCluster</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> )</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>操作</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>redis</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>集群</w:t></w:r>' + `
  '</w:p>'
$p2.InsertXML($xml2)

# --- Paragraph 3: "【参考】" + newly cited source files + bookmark ----------
$p3 = $d.Paragraphs(3).Range
$xml3 = '<w:p xmlns:w="' + $wNs + '" w:rsidR="00B72A5A" w:rsidRDefault="00B72A5A">' + `
  '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:noProof/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:noProof/></w:rPr><w:t>【参考】</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:noProof/></w:rPr><w:t>jedis/RedisTest.java</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:noProof/></w:rPr><w:t>、</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:noProof/></w:rPr><w:t>jedis/JedisPubSubList</w:t></w:r>' + `
  '<w:r><w:rPr><w:noProof/></w:rPr><w:t>ener.java</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$p3.InsertXML($xml3)

# --- Append a new empty paragraph right after paragraph 3, before sectPr ---
$p3end = $d.Paragraphs(3).Range.End
$tail = $d.Range($p3end, $p3end)
$tail.InsertXML('<w:p xmlns:w="' + $wNs + '"/>')
